$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-46) holds the "Förändrad" (changed) date, stored as
# serial date 45189 (2023-09-20). Bump it by one day to 45190 (2023-09-21)
# for every data row, matching the diff.
for ($row = 2; $row -le 46; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
